$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New normalization pairs appended to the word-mapping table
# (col A = informal/abbreviated form, col B = normalized form)
$ws.Range("A40").Value = "dah"
$ws.Range("B40").Value = "sudah"

$ws.Range("A41").Value = "nak"
$ws.Range("B41").Value = "akan"

$ws.Range("A42").Value = "tu"
$ws.Range("B42").Value = "itu"

# Match the author's final scroll position / selection
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B42").Select()
